# Add a new "MonkeyBusiness" worksheet (test table) after "SkullIsland",
# matching the commit "Other: add test table for MonkeyBussiness Theme".

$wb = $excel.ActiveWorkbook

# --- Leave the "SkullIsland" sheet the way it was found, just record it so
#     we can select its first row before we switch focus away from it (this
#     matches the row-1 selection state captured in the target file). ---
$skull = $wb.Worksheets.Item("SkullIsland")
$skull.Activate()
$skull.Rows.Item(1).Select()

# --- Create the new sheet as the last tab in the workbook. ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "MonkeyBusiness"

# --- Header row ---
$ws.Range("A1").Value = "Skull"
$ws.Range("B1").Value = "Monkey"
$ws.Range("C1").Value = "Parrot"
$ws.Range("D1").Value = "Sword"
$ws.Range("E1").Value = "Coin"
$ws.Range("F1").Value = "Diamond"
$ws.Range("G1").Value = "expect"
$ws.Range("H1").Value = "skull from card"
$ws.Range("I1").Value = "sum check"
$ws.Range("J1").Value = "note"

# --- Data rows ---
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 200
$ws.Range("J2").Value = "coinx2"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 1
$ws.Range("G3").Value = 200
$ws.Range("J3").Value = "3oak+coin"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("G4").Value = 200
$ws.Range("J4").Value = "3oak+coin"

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("G5").Value = 200
$ws.Range("J5").Value = "4oak"

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 1
$ws.Range("G6").Value = 500
$ws.Range("J6").Value = "5oak"

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 3
$ws.Range("G7").Value = 1000
$ws.Range("J7").Value = "6oak"

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 6
$ws.Range("G8").Value = 2000
$ws.Range("J8").Value = "7oak"

$ws.Range("A9").Value = 0
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = 5
$ws.Range("G9").Value = 5500
$ws.Range("J9").Value = "8oak+FC"

$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 3
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 1300
$ws.Range("J10").Value = "5oak+coinx3+FC"

# --- "sum check" formulas, column I ---
$ws.Range("I2").Formula = "=SUM(A2:F2)"
$ws.Range("I3").Formula = "=SUM(A3:F3)"
$ws.Range("I4").Formula = "=SUM(A4:F4)"
$ws.Range("I5").Formula = "=SUM(A5:F5)"
$ws.Range("I6").Formula = "=SUM(A6:F6)"
$ws.Range("I7").Formula = "=SUM(A7:F7)"
$ws.Range("I8").Formula = "=SUM(A8:F8)"
$ws.Range("I9").Formula = "=SUM(A9:F9)"
$ws.Range("I10").Formula = "=SUM(A10:F10)"

# --- Column widths (best-fit on the "skull from card" / "note" columns) ---
$ws.Columns.Item(8).ColumnWidth = 13.33203125
$ws.Columns.Item(10).ColumnWidth = 15.5

# --- Selection state left on the new sheet, and make it the active tab. ---
$ws.Activate()
$ws.Range("I11:I14").Select()
